$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the header row (row 2) with clarified column descriptions.
#    The new shared-string entries must be created in this exact order so
#    that they land on shared-string indices 23..27 (Ik (ipsc LPE), ZS (Z LPE),
#    ZL (Z LN), Ik (Ipsc LN), uln (V)).
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "Ik (A) (ipsc (LPE))"
$ws.Range("E2").Value = "ZS (Ω) (Z (LPE))"
$ws.Range("B2").Value = "ZL (Ω) (Z (LN))"
$ws.Range("C2").Value = "Ik(A) (Ipsc (LN))"
$ws.Range("A2").Value = "uln (V)"

# ---------------------------------------------------------------------------
# 2) New block of rows describing "AUTO TN je vsak zase" with a second copy
#    of the header (shifted one column to the right, rows 19-21) bracketed
#    by two explanatory labels (rows 18 and 24).
# ---------------------------------------------------------------------------

# Row 18: section label
$ws.Range("A18").Value = "del poti npr. razdelilec 2"

# Row 19: label + full header copy (columns B..O mirror columns A..N of row 2)
$ws.Range("A19").Value = "del poti npr. F32"
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4163) | Out-Null   # xlPasteValues

# Row 20: identical repeat of row 19
$ws.Range("A20").Value = "del poti npr. F32"
$ws.Range("A19:O19").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("B19:O19").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4163) | Out-Null

# Row 21: identical repeat of row 19/20
$ws.Range("A21").Value = "del poti npr. F32"
$ws.Range("A19:O19").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("B19:O19").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4163) | Out-Null

# Row 24: closing note
$ws.Range("A24").Value = "AUTO TN je vsak zase, skupaj z največ 1 r low 4"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Column widths for the newly used columns (best achievable approximation
#    given the automation layer's internal rounding of ColumnWidth).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.83
$ws.Columns.Item(3).ColumnWidth = 15.5
$ws.Columns.Item(5).ColumnWidth = 11.67
$ws.Columns.Item(6).ColumnWidth = 13.17
$ws.Columns.Item(8).ColumnWidth = 20.33

# ---------------------------------------------------------------------------
# 4) Restore the selection to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("D26").Select() | Out-Null
